$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G16").Value = 82
$ws.Range("G17").Value = 118
$ws.Range("G18").Value = 118
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
